# Applies the "Clean up code, fix authentication, update to use recipes.csv,
# remove unnecessary console logs" data refresh to the Users sheet:
#  - remove the hyperlink that was attached to B3 (sanjanaplayz@gmail.com row)
#  - delete the two trailing rows (old Test22 / Test3 seed rows)
#  - replace the remaining seed data with the new test/test1/test2 users
#  - clear the leftover hyperlink font/style from B3
#  - move the active selection the way the source file has it

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the mailto: hyperlink (and its formatting) that lived on B3.
$ws.Hyperlinks.Delete()
$ws.Range("B3").ClearFormats()

# Remove the last two seed rows (ids 4 and 5) entirely.
$ws.Rows("5:6").Delete()

# --- Row 2 (id 1) ---------------------------------------------------
$ws.Range("B2").Value = "test@example.com"
$ws.Range("C2").Value = "`$2b`$10`$2eitSKhyJFA.BJYmgUkz8un0N5KX7U7JF0wb/IlBAkvugAWy7OGgG"
$ws.Range("D2").Value = "testuser"
$ws.Range("E2").Value = 45998.781270277774

# --- Row 3 (id 2) ---------------------------------------------------
$ws.Range("B3").Value = "test1@gmail.com"
$ws.Range("C3").Value = "`$2b`$10`$jQbw0KhcDELathahjo/nzOvZRwIZvXLQCn6Tgn19EN2ln83ybghsC"
$ws.Range("D3").Value = "test1"
$ws.Range("E3").Value = 45998.781803784717

# --- Row 4 (id 3) ----------------------------------------------------
$ws.Range("B4").Value = "test2@gmail.com"
$ws.Range("C4").Value = "`$2b`$10`$nS2qkMGkTSADbFDEP6AQkuEpZVBf4VZ4HLGB4g/A4QYygoL19rZ8e"
$ws.Range("D4").Value = "test2"
$ws.Range("E4").Value = 45998.784592118056

# Match the saved file's view state (selection on E10, tab selected).
$ws.Range("E10").Select()
